$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 34

# Duplicate the previous (last) data row so formatting and the
# "empty" F/H cells match exactly, then overwrite the cells whose
# values actually changed for this new run.
$srcRange = $ws.Range("A33:H33")
$dstRange = $ws.Range("A34:H34")
$srcRange.Copy($dstRange)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = "2025-08-20 03:51:22 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-20 09:21:22 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($row, 7).Value = 0

$wb.Save()
